# Insert a new weekly data row at row 177 (pushing existing rows 177-262 down to 178-263).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("177:177").Insert()

$ws.Range("A177").Value = 10
$ws.Range("B177").Value = "Vega Modelo de Temuco"
$ws.Range("C177").Value = "La Araucanía"
$ws.Range("D177").NumberFormat = $ws.Range("D178").NumberFormat
$ws.Range("D177").Value = 44518
$ws.Range("E177").Value = 9
$ws.Range("F177").Value = 100112040
$ws.Range("G177").Value = "Cilantro"
$ws.Range("H177").Value = "Sin especificar"
$ws.Range("I177").Value = "Primera"
$ws.Range("J177").Value = 85
$ws.Range("K177").Value = 5000
$ws.Range("L177").Value = 6000
$ws.Range("M177").Value = 5412
$ws.Range("N177").Value = "$/docena de atados (2 kilos)"
$ws.Range("O177").Value = "Provincia de Cautín"
$ws.Range("P177").Value = 2706
$ws.Range("Q177").Value = 2
$ws.Range("R177").Value = "Hortaliza"
